$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update username value in A2 from "qatitans65" to "qatitans78"
$ws.Range("A2").Value = "qatitans78"

# Move the selection from D1 to A2
$ws.Range("A2").Select()
